$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.175.14"
$ws.Range("E2").Value = "  -1.14%  "
$ws.Range("D3").Value = "2.434.90"
$ws.Range("E3").Value = "  -1.56%  "
$ws.Range("E4").Value = "  +0.13%  "
$c = $ws.Range("D5")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "316.51"
$c.Style = $origStyle
$ws.Range("E5").Value = "  -0.72%  "
$c = $ws.Range("D6")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "89.10"
$c.Style = $origStyle
$ws.Range("E6").Value = "  -3.96%  "
$ws.Range("E9").Value = "  -4.09%  "
$c = $ws.Range("D10")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "32.15"
$c.Style = $origStyle
$ws.Range("E10").Value = "  -2.68%  "
$ws.Range("E11").Value = "  -4.35%  "
$ws.Range("E12").Value = "  -2.70%  "
$ws.Range("D13").Value = "2.809.68"
$ws.Range("E13").Value = "  -1.60%  "
$ws.Range("E14").Value = "  -2.93%  "
$c = $ws.Range("D15")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "15.59"
$c.Style = $origStyle
$ws.Range("E15").Value = "  -0.42%  "
$ws.Range("D16").Value = "2.441.97"
$ws.Range("E16").Value = "  -1.74%  "
$c = $ws.Range("D17")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.773"
$c.Style = $origStyle
$ws.Range("E17").Value = "  -2.28%  "
$ws.Range("D18").Value = "41.116.77"
$ws.Range("E18").Value = "  -1.20%  "
$ws.Range("D19").Value = "0.0₃0924"
$ws.Range("E19").Value = "  -3.64%  "
$c = $ws.Range("D20")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "6.25"
$c.Style = $origStyle
$ws.Range("E20").Value = "  -3.90%  "
$c = $ws.Range("D21")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "72.17"
$c.Style = $origStyle
$ws.Range("E21").Value = "  +1.08%  "
$c = $ws.Range("D22")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "10.99"
$c.Style = $origStyle
$ws.Range("E22").Value = "  -4.37%  "
$c = $ws.Range("D23")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "235.22"
$c.Style = $origStyle
$ws.Range("E23").Value = "  -2.63%  "
$c = $ws.Range("D24")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.69"
$c.Style = $origStyle
$ws.Range("E24").Value = "  -2.53%  "
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("E26").Value = "  -2.77%  "
$c = $ws.Range("D27")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "24.03"
$c.Style = $origStyle
$ws.Range("E27").Value = "  -3.25%  "
$c = $ws.Range("D28")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.22"
$c.Style = $origStyle
$ws.Range("E28").Value = "  -3.13%  "
$ws.Range("E29").Value = "  -3.61%  "
$c = $ws.Range("D30")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "34.85"
$c.Style = $origStyle
$ws.Range("E30").Value = "  -4.64%  "
$c = $ws.Range("D31")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "156.39"
$c.Style = $origStyle
$ws.Range("E31").Value = "  -0.41%  "
$ws.Range("E32").Value = "  +0.01%  "
$c = $ws.Range("D33")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "5.27"
$c.Style = $origStyle
$ws.Range("E33").Value = "  -4.81%  "
$ws.Range("E34").Value = "  -2.28%  "
$c = $ws.Range("D35")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.0744"
$c.Style = $origStyle
$ws.Range("E35").Value = "  -3.66%  "
$c = $ws.Range("D36")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.93"
$c.Style = $origStyle
$ws.Range("E36").Value = "  +0.25%  "
$c = $ws.Range("D37")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "16.66"
$c.Style = $origStyle
$ws.Range("E37").Value = "  -4.54%  "
$ws.Range("E38").Value = "  -0.74%  "
$c = $ws.Range("D39")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.78"
$c.Style = $origStyle
$ws.Range("E39").Value = "  -3.39%  "
$c = $ws.Range("D40")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.100"
$c.Style = $origStyle
$ws.Range("E40").Value = "  -3.19%  "
$ws.Range("E41").Value = "  -3.94%  "
$ws.Range("E42").Value = "  -6.74%  "
$ws.Range("D43").Value = "1.987.29"
$ws.Range("E43").Value = "  +0.13%  "
$c = $ws.Range("D44")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.0275"
$c.Style = $origStyle
$ws.Range("E44").Value = "  -3.64%  "
$c = $ws.Range("D45")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "18.39"
$c.Style = $origStyle
$ws.Range("E45").Value = "  -4.50%  "
$c = $ws.Range("D46")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.88"
$c.Style = $origStyle
$ws.Range("E46").Value = "  -5.00%  "
$c = $ws.Range("D47")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "9.45"
$c.Style = $origStyle
$ws.Range("E47").Value = "  +2.33%  "
$ws.Range("D48").Value = "2.668.72"
$ws.Range("E48").Value = "  -1.65%  "
$c = $ws.Range("D49")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "95.24"
$c.Style = $origStyle
$ws.Range("E49").Value = "  -2.47%  "
$c = $ws.Range("D50")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "73.18"
$c.Style = $origStyle
$ws.Range("E50").Value = "  -0.94%  "
$c = $ws.Range("D51")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "51.49"
$c.Style = $origStyle
$ws.Range("E51").Value = "  -2.26%  "
